# The "GitHub Link:" paragraph originally reads, across several runs:
#   "GitHub Link: " + "https:" + " " + "//github.com/" + "okunadeaminat" + "-ai/my-smart-daily-planner"
#
# The edit collapses the last three text fragments (the "//github.com/",
# the spell-check-flagged "okunadeaminat", and the "-ai/my-smart-daily-planner"
# runs) into a single run whose text is the complete URL
# "https://github.com/okunadeaminat-ai/my-smart-daily-planner", leaving the
# "GitHub Link: https: " prefix untouched.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "//github.com/okunadeaminat-ai/my-smart-daily-planner",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://github.com/okunadeaminat-ai/my-smart-daily-planner",
    2)

Write-Output "GitHub link replaced: $found"

if (-not $found) {
    throw "Could not find the GitHub link text to replace."
}
